$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = "Datos actualizados a 18 de Julio de 2020 a las 14:01"

# Row 4
$ws.Range("B4").Value = 3771101
$ws.Range("C4").Value = 1089
$ws.Range("D4").Value = 1741626
$ws.Range("E4").Value = 1887395
$ws.Range("G4").Value = 16
$ws.Range("H4").Value = 142080

# Row 5
$ws.Range("B5").Value = 2049140
$ws.Range("C5").Value = 443
$ws.Range("E5").Value = 604401
$ws.Range("G5").Value = 32
$ws.Range("H5").Value = 77964

# Row 6
$ws.Range("B6").Value = 1044963
$ws.Range("C6").Value = 4506
$ws.Range("D6").Value = 655667
$ws.Range("E6").Value = 362949
$ws.Range("G6").Value = 62
$ws.Range("H6").Value = 26347

# Row 19
$ws.Range("B19").Value = 202372
$ws.Range("C19").Value = 27
$ws.Range("E19").Value = 5712

# Row 33
$ws.Range("B33").Value = 65953
$ws.Range("C33").Value = 171
$ws.Range("D33").Value = 57856
$ws.Range("E33").Value = 7602
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = 495

# Row 37
$ws.Range("B37").Value = 58904
$ws.Range("C37").Value = 683
$ws.Range("D37").Value = 49020
$ws.Range("E37").Value = 9477
$ws.Range("G37").Value = 3
$ws.Range("H37").Value = 407

# Row 54
$ws.Range("B54").Value = 33492
$ws.Range("C54").Value = 110
$ws.Range("E54").Value = 1623

# Row 64
$ws.Range("B64").Value = 17502
$ws.Range("C64").Value = 57
$ws.Range("D64").Value = 11637
$ws.Range("E64").Value = 5825

# Row 73
$ws.Range("B73").Value = 12750
$ws.Range("C73").Value = 688
$ws.Range("D73").Value = 4440
$ws.Range("E73").Value = 8085
$ws.Range("G73").Value = 3
$ws.Range("H73").Value = 225

# Row 75
$ws.Range("B75").Value = 11441
$ws.Range("C75").Value = 206
$ws.Range("D75").Value = 8161

# Row 77
$ws.Range("B77").Value = 10682
$ws.Range("C77").Value = 83
$ws.Range("D77").Value = 5637
$ws.Range("E77").Value = 4372
$ws.Range("G77").Value = 2
$ws.Range("H77").Value = 673

# Row 83
$ws.Range("B83").Value = 8669
$ws.Range("C83").Value = 125
$ws.Range("D83").Value = 5859
$ws.Range("E83").Value = 2647
$ws.Range("G83").Value = 3
$ws.Range("H83").Value = 163

# Row 86
$ws.Range("A86").Value = "Estado de Palestina"
$ws.Range("B86").Value = 8200
$ws.Range("C86").Value = 436
$ws.Range("D86").Value = 1596
$ws.Range("E86").Value = 6549
$ws.Range("G86").Value = 2
$ws.Range("H86").Value = 55

# Row 87
$ws.Range("A87").Value = "Bosnia y Herzegovina"
$ws.Range("B87").Value = 8161
$ws.Range("C87").Value = 253
$ws.Range("D87").Value = 3648
$ws.Range("E87").Value = 4267
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 246

# Row 90
$ws.Range("A90").Value = "Madagascar"
$ws.Range("B90").Value = 6849
$ws.Range("C90").Value = 382
$ws.Range("D90").Value = 3339
$ws.Range("E90").Value = 3455
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 55

# Row 91
$ws.Range("A91").Value = "Tayikistan"
$ws.Range("B91").Value = 6786
$ws.Range("D91").Value = 5483
$ws.Range("E91").Value = 1247
$ws.Range("H91").Value = 56

# Row 92
$ws.Range("A92").Value = "Guayana Francesa"
$ws.Range("B92").Value = 6509
$ws.Range("D92").Value = 3932
$ws.Range("E92").Value = 2543
$ws.Range("H92").Value = 34

# Row 109
$ws.Range("A109").Value = "Malaui"
$ws.Range("C109").Value = 5
$ws.Range("D109").Value = 1111
$ws.Range("E109").Value = 1644
$ws.Range("H109").Value = 55

# Row 110
$ws.Range("A110").Value = "Zambia"
$ws.Range("B110").Value = 2810
$ws.Range("D110").Value = 1450
$ws.Range("E110").Value = 1251
$ws.Range("H110").Value = 109

# Row 112
$ws.Range("A112").Value = "Sri Lanka"
$ws.Range("B112").Value = 2701
$ws.Range("C112").Value = 4
$ws.Range("D112").Value = 2023
$ws.Range("E112").Value = 667
$ws.Range("H112").Value = 11

# Row 113
$ws.Range("A113").Value = "Libano"
$ws.Range("B113").Value = 2700
$ws.Range("D113").Value = 1485
$ws.Range("E113").Value = 1175
$ws.Range("H113").Value = 40

# Row 124
$ws.Range("B124").Value = 1922
$ws.Range("C124").Value = 6
$ws.Range("D124").Value = 1902
$ws.Range("E124").Value = 10

# Row 125
$ws.Range("B125").Value = 1915
$ws.Range("C125").Value = 7
$ws.Range("D125").Value = 1600
$ws.Range("E125").Value = 235
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 80

# Row 142
$ws.Range("B142").Value = 1062
$ws.Range("C142").Value = 6
$ws.Range("E142").Value = 39

# Row 143
$ws.Range("B143").Value = 1047
$ws.Range("C143").Value = 2
$ws.Range("E143").Value = 107

# Row 155
$ws.Range("B155").Value = 675
$ws.Range("C155").Value = 1
$ws.Range("E155").Value = 4

# Row 163
$ws.Range("A163").Value = "Lesoto"
$ws.Range("B163").Value = 359
$ws.Range("C163").Value = 48
$ws.Range("D163").Value = 69
$ws.Range("E163").Value = 284
$ws.Range("H163").Value = 6

# Row 164
$ws.Range("A164").Value = "Mauricio"
$ws.Range("B164").Value = 343
$ws.Range("D164").Value = 331
$ws.Range("E164").Value = 2
$ws.Range("H164").Value = 10

# Row 165
$ws.Range("A165").Value = "Birmania"
$ws.Range("B165").Value = 339
$ws.Range("D165").Value = 271
$ws.Range("E165").Value = 62
$ws.Range("H165").Value = 6

# Row 166
$ws.Range("A166").Value = "Isla de Man"
$ws.Range("B166").Value = 336
$ws.Range("D166").Value = 312
$ws.Range("E166").Value = 0
$ws.Range("H166").Value = 24

# Row 167
$ws.Range("A167").Value = "Comoras"
$ws.Range("B167").Value = 328
$ws.Range("D167").Value = 311
$ws.Range("E167").Value = 10
$ws.Range("H167").Value = 7

# Row 168
$ws.Range("A168").Value = "Guyana"
$ws.Range("B168").Value = 320
$ws.Range("D168").Value = 156
$ws.Range("E168").Value = 145
$ws.Range("H168").Value = 19

# Row 210
$ws.Range("A210").Value = "Islas Malvinas"

# Row 211
$ws.Range("A211").Value = "Groenlandia"
